# Rename the category label "congenital" to "misc_long_term" across all
# worksheets in the workbook. The label lives in either cell A3 or A4
# (depending on the sheet's layout), so scan both candidate cells on every
# worksheet and update any that match.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in @("A3", "A4")) {
        $cell = $ws.Range($addr)
        if ($cell.Text -eq "congenital") {
            $cell.Value = "misc_long_term"
        }
    }
}
